# Applies the numeric updates to column F ("reads"/popularity counters)
# across the 展览, 演出 and 全部类型 sheets, as captured by the source diff.

$wb = $excel.ActiveWorkbook

# Map of sheet name -> list of (CellAddress, NewValue)
$changes = @{
    "展览" = @(
        @{ Cell = "F7";  Value = 447 },
        @{ Cell = "F12"; Value = 1047 },
        @{ Cell = "F16"; Value = 1492 },
        @{ Cell = "F18"; Value = 226 },
        @{ Cell = "F21"; Value = 812 },
        @{ Cell = "F22"; Value = 1150 },
        @{ Cell = "F24"; Value = 1916 },
        @{ Cell = "F25"; Value = 2653 },
        @{ Cell = "F28"; Value = 34 },
        @{ Cell = "F31"; Value = 1213 },
        @{ Cell = "F33"; Value = 1342 },
        @{ Cell = "F38"; Value = 667 },
        @{ Cell = "F39"; Value = 835 },
        @{ Cell = "F41"; Value = 244 }
    )
    "演出" = @(
        @{ Cell = "F23"; Value = 15 }
    )
    "全部类型" = @(
        @{ Cell = "F11"; Value = 447 },
        @{ Cell = "F21"; Value = 1492 },
        @{ Cell = "F23"; Value = 226 },
        @{ Cell = "F26"; Value = 1150 },
        @{ Cell = "F27"; Value = 2653 },
        @{ Cell = "F32"; Value = 34 },
        @{ Cell = "F35"; Value = 0 },
        @{ Cell = "F36"; Value = 1213 },
        @{ Cell = "F40"; Value = 1342 },
        @{ Cell = "F43"; Value = 667 },
        @{ Cell = "F44"; Value = 835 },
        @{ Cell = "F47"; Value = 15 },
        @{ Cell = "F48"; Value = 244 }
    )
}

foreach ($sheetName in $changes.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($change in $changes[$sheetName]) {
        $ws.Range($change.Cell).Value = $change.Value
    }
}
